# Updates cryptos list values (price & volume/1h columns) and
# replaces a few coin rows (48, 49, 51), matching the
# 'Updated cryptos list ... with GitHub Actions' commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.793.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -6.38%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.937.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -8.82%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.94%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -11.77%  "

# Row 7
$ws.Range("E7").Value = "  -0.33%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.916.03"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -9.45%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.466"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -14.58%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.146"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -16.45%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.87"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -10.46%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.436"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -12.46%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -16.83%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000204"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -16.84%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.399.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -9.17%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.681.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.59%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.109"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.39%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.932.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -8.87%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "469.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -12.30%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -13.32%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -13.31%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.637"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -16.51%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -16.84%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -11.72%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -13.89%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -18.67%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -13.68%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.54%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "24.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -15.73%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -10.21%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.997"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.64%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "480.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -12.38%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.02%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -15.48%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -15.96%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0389"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.12%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.116"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.65%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0752"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -13.13%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -16.05%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.712.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.24%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.16%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -12.54%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.225"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -15.00%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "112.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.63%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.102"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.35%  "

# Row 48 (was Fetch.AI -> now PEPE)
$ws.Range("B48").Value = "PEPE"
$ws.Range("C48").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₃0484"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -18.04%  "

# Row 49 (was PEPE -> now Fetch.AI)
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -13.69%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -17.38%  "

# Row 51 (was BitgetToken -> now ThetaToken)
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -20.87%  "
